$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/16/2024  Through  12/22/2024"

# --- Donor cells carrying the text-placeholder style (style 13) ---
$donorZero = "C14"
$donorNA = "E14"

# Row 14
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "0"
$ws.Range($donorZero).Copy() | Out-Null
$ws.Range("F14").PasteSpecial(-4122)

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range($donorZero).Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range($donorNA).Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("N15").Value = -75

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -44.444444444444
$ws.Range("I16").Value = 184
$ws.Range("J16").Value = 208
$ws.Range("K16").Value = -11.538461538461
$ws.Range("L16").Value = -1.604278074866
$ws.Range("M16").Value = -56.398104265402
$ws.Range("N16").Value = -91.006842619745

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -70
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = -27.586206896551
$ws.Range("I17").Value = 417
$ws.Range("J17").Value = 412
$ws.Range("K17").Value = 1.213592233009
$ws.Range("L17").Value = -6.081081081081
$ws.Range("M17").Value = -0.23923444976
$ws.Range("N17").Value = -62.466246624662

# Row 18
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 11
$ws.Range("H18").Value = -8.333333333333
$ws.Range("I18").Value = 195
$ws.Range("J18").Value = 226
$ws.Range("K18").Value = -13.716814159292
$ws.Range("L18").Value = -19.753086419753
$ws.Range("M18").Value = -51.851851851851
$ws.Range("N18").Value = -84.261501210653

# Row 19
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -44.444444444444
$ws.Range("F19").Value = 22
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = -24.137931034482
$ws.Range("I19").Value = 362
$ws.Range("J19").Value = 436
$ws.Range("K19").Value = -16.97247706422
$ws.Range("L19").Value = -11.274509803921
$ws.Range("M19").Value = -9.950248756218
$ws.Range("N19").Value = -57.511737089201

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -60
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = -86.363636363636
$ws.Range("I20").Value = 104
$ws.Range("J20").Value = 137
$ws.Range("K20").Value = -24.087591240875
$ws.Range("L20").Value = 4
$ws.Range("M20").Value = -9.565217391304
$ws.Range("N20").Value = -80.847145488029

# Row 21
$ws.Range("C21").Value = 15
$ws.Range("E21").Value = -46.428571428571
$ws.Range("F21").Value = 68
$ws.Range("G21").Value = 111
$ws.Range("H21").Value = -38.738738738738
$ws.Range("I21").Value = 1289
$ws.Range("J21").Value = 1451
$ws.Range("K21").Value = -11.164713990351
$ws.Range("L21").Value = -8.904593639575
$ws.Range("M21").Value = -28.508042151968
$ws.Range("N21").Value = -78.226351351351

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range($donorZero).Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range($donorNA).Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("M22").Value = -57.575757575757

# Row 23
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -83.333333333333
$ws.Range("F23").Value = 16
$ws.Range("G23").Value = 22
$ws.Range("H23").Value = -27.272727272727
$ws.Range("I23").Value = 206
$ws.Range("J23").Value = 245
$ws.Range("K23").Value = -15.918367346938
$ws.Range("L23").Value = -8.035714285714
$ws.Range("M23").Value = -2.830188679245

# Row 24
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 47.619047619047
$ws.Range("F24").Value = 107
$ws.Range("G24").Value = 95
$ws.Range("H24").Value = 12.631578947368
$ws.Range("I24").Value = 1442
$ws.Range("J24").Value = 1425
$ws.Range("K24").Value = 1.19298245614
$ws.Range("L24").Value = -9.308176100628
$ws.Range("M24").Value = 41.929133858267

# Row 25
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 266.666666666667
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = 33.333333333333
$ws.Range("I25").Value = 653
$ws.Range("J25").Value = 521
$ws.Range("K25").Value = 25.335892514395
$ws.Range("L25").Value = -5.772005772005

# Row 26
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = -41.666666666666
$ws.Range("F26").Value = 47
$ws.Range("G26").Value = 45
$ws.Range("H26").Value = 4.444444444444
$ws.Range("I26").Value = 610
$ws.Range("J26").Value = 626
$ws.Range("K26").Value = -2.55591054313
$ws.Range("L26").Value = -4.984423676012
$ws.Range("M26").Value = -35.58606124604

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range($donorZero).Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range($donorNA).Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122)

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("F28").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = 3
$ws.Range("F28").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = -66.666666666666
$ws.Range("H28").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -25
$ws.Range("I28").Value = 50
$ws.Range("J28").Value = 45
$ws.Range("K28").Value = 11.111111111111
$ws.Range("L28").Value = 16.279069767441

# Row 29
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range($donorZero).Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("F29").Value = 3
$ws.Range("H29").Value = 50
$ws.Range("L29").Value = -14.705882352941
$ws.Range("N29").Value = -89.338235294117

# Row 30
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range($donorZero).Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("F30").Value = 2
$ws.Range("H30").Value = 0
$ws.Range("L30").Value = -14.814814814814
$ws.Range("N30").Value = -90.416666666666

# Row 31
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = -50
$ws.Range("J31").Value = 11
$ws.Range("K31").Value = 36.363636363636

# --- Insert new blank row at 56 (shifts old 56/57 footer down to 57/58) ---
$ws.Rows(56).Insert()
$ws.Range("A56").Clear()

